# Generate Report for Handback
# Update the timestamp values on the "Overview", "zh-cn" and "de-de" sheets
# to reflect the latest handoff/handback generation times.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file row
$wsOverview.Range("G2").Value = "2016-08-24 05:03:11"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the first file row
$wsZhCn.Range("H2").Value = "2016-08-24 05:03:03"
$wsZhCn.Range("K2").Value = "2016-08-24 05:03:29"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the first file row
$wsDeDe.Range("H2").Value = "2016-08-24 05:03:11"
$wsDeDe.Range("K2").Value = "2016-08-24 05:03:36"
